$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for B2:K21 (one row per dictionary entry, keyed by row number)
$data = @{
    2 = @(-19.14671036488678,4.321925733515622,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    3 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,1.214544862311391,-19.14671036488678,-19.14671036488678)
    4 = @(-19.14671036488678,-19.14671036488678,1.133023757235197,-19.14671036488678,3.489500248710589,-19.14671036488678,1.521541868447137,-19.14671036488678,0.8765255458390527,-19.14671036488678)
    5 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,2.826879457197482,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    6 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    7 = @(2.57438757185137,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    8 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,1.523843023544502,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    9 = @(3.811860226949741,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    10 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,1.71693525441244,-19.14671036488678,2.217411997093501)
    11 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,2.911844602411979,-19.14671036488678,2.837100533465396,-19.14671036488678,-19.14671036488678,-19.14671036488678,1.954296454649256)
    12 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    13 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,2.599564885765298,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,1.646242727886163,1.769338915681817)
    14 = @(-19.14671036488678,-19.14671036488678,1.678267792295479,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,1.954322110987776)
    15 = @(-19.14671036488678,-19.14671036488678,1.869440249378016,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678)
    16 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,1.941436683644,-19.14671036488678)
    17 = @(-19.14671036488678,-19.14671036488678,1.979286342985898,-19.14671036488678,-19.14671036488678,-19.14671036488678,2.04997159806794,2.132422386673785,2.525129982595545,-19.14671036488678)
    18 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,-19.14671036488678,2.0141267917423,2.063628634727881,2.442736428999218,-19.14671036488678)
    19 = @(-19.14671036488678,-19.14671036488678,2.12586531892233,-19.14671036488678,-19.14671036488678,-19.14671036488678,1.644828547531502,1.817904778658912,-19.14671036488678,-19.14671036488678)
    20 = @(-19.14671036488678,-19.14671036488678,1.403172987208867,-19.14671036488678,3.132297376687202,-19.14671036488678,1.628826777538171,1.202476899411775,-19.14671036488678,2.066883378918029)
    21 = @(-19.14671036488678,-19.14671036488678,-19.14671036488678,1.822761178741657,-19.14671036488678,2.525760298859639,1.44834036859698,-19.14671036488678,-19.14671036488678,-19.14671036488678)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $rowVals[$i]
    }
}
